$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") bumped by one day (46072 -> 46073) for all data rows
$ws.Range("C2:C9").Value = 46073

# Rows 3 and 4 swap their A/B/G values
$ws.Range("A3").Value = "A 25617-2024"
$ws.Range("B3").Value = 45463
$ws.Range("G3").Value = 2.3

$ws.Range("A4").Value = "A 34310-2024"
$ws.Range("B4").Value = 45524
$ws.Range("G4").Value = 4.8

# Rows 5 and 6 swap their A/B/G values
$ws.Range("A5").Value = "A 54782-2022"
$ws.Range("B5").Value = 44883
$ws.Range("G5").Value = 5.5

$ws.Range("A6").Value = "A 45983-2023"
$ws.Range("B6").Value = 45196
$ws.Range("G6").Value = 0.6
